$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2026-02 (row 27)
$ws.Range("B27").Value = 6567
$ws.Range("C27").Value = 1018
$ws.Range("D27").Value = 6127524
$ws.Range("E27").Value = 933.0781178620375
$ws.Range("F27").Value = 10.36974789915965
$ws.Range("G27").Value = 7.4973600844773
$ws.Range("H27").Value = 25.67972930815194
